$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.098.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.21%  "

$ws.Range("D3").Value = "'2.473.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.23%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'573.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").Value = "'148.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "'0.531"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.60%  "

$ws.Range("E9").Value = "  -0.19%  "

$ws.Range("D10").Value = "'0.163"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.23%  "

$ws.Range("D11").Value = "'5.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.28%  "

$ws.Range("D12").Value = "'0.348"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.34%  "

$ws.Range("D13").Value = "'28.96"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("D14").Value = "'0.0000176"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.91%  "

$ws.Range("D15").Value = "'2.911.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("D16").Value = "'62.945.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").Value = "'2.463.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").Value = "'7.70"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.14%  "

$ws.Range("D19").Value = "'10.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.25%  "

$ws.Range("D20").Value = "'2.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.02%  "

$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").Value = "'4.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "'322.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.15%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'10.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.98%  "

$ws.Range("D25").Value = "'65.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("D26").Value = "'653.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.91%  "

$ws.Range("D27").Value = "'2.587.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.26%  "

$ws.Range("D28").Value = "'0.0₃0974"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.42%  "

$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").Value = "'1.43"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.54%  "

$ws.Range("D31").Value = "'7.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.16%  "

$ws.Range("D32").Value = "'1.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.41%  "

$ws.Range("E33").Value = "  +0.18%  "

$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").Value = "'1.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.70%  "

$ws.Range("D36").Value = "'4.69"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.00%  "

$ws.Range("D37").Value = "'5.40"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.22%  "

$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "'0.366"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.73%  "

$ws.Range("B39").Value = "EthereumClassic"
$ws.Range("C39").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'18.58"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.14%  "

$ws.Range("D40").Value = "'149.55"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.38%  "

$ws.Range("D41").Value = "'2.72"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.26%  "

$ws.Range("D42").Value = "'1.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.69%  "

$ws.Range("D43").Value = "'0.0₆0311"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.02%  "

$ws.Range("E44").Value = "  -0.02%  "

$ws.Range("D45").Value = "'153.60"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("D46").Value = "'15.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.87%  "

$ws.Range("D47").Value = "'3.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.16%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'20.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.02%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.606"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").Value = "'0.0510"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.77%  "

$ws.Range("D51").Value = "'0.0906"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.44%  "
